$d = $word.ActiveDocument

# --- 1. Insert a new "Meta description" paragraph right after the title (Heading1) ---
# First split off a new (blank) paragraph right after the title, then overwrite that
# blank paragraph's OOXML with our target run structure. Supplying a <w:p> with no
# <w:pPr> clears the inherited Heading1 style cleanly (no explicit pStyle left behind,
# matching the plain body paragraphs elsewhere in the document) without leaving stray
# rsid bookkeeping attributes the way assigning .Style afterwards would.
$titlePara = $d.Paragraphs.Item(1)
$newPara = $titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)

$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Discover Dragon''s Inferno online slot by WMS with free play, RTP, and features. Read our unbiased review of 2021.</w:t></w:r></w:p>'
$metaPara.Range.InsertXML($metaXml)

# --- 2. Remove the duplicated bold title paragraph near the end of the document ---
$dupTitlePara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$dupTitlePara.Range.Delete()

# --- 3. Replace the meta-description text in the final (italic) paragraph with the new
#        image-prompt text, editing only the text run so the leading empty run and the
#        italic formatting survive untouched. Using the run's own Range (not Find/Replace)
#        also avoids Word's smart-quote autocorrection mangling the apostrophes. ---
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$oldText = "Discover Dragon's Inferno online slot by WMS with free play, RTP, and features. Read our unbiased review of 2021."
$newText = "Create a feature image for Dragon's Inferno that features a happy Maya warrior wearing glasses in a cartoon style. The warrior should be holding a dragon's egg in one hand and a treasure chest overflowing with gold coins in the other. In the background, there should be a fiery landscape with dragons flying in the distance. The image should be bright and colorful, capturing the adventurous spirit and fun gameplay of this slot machine game."

$paraEnd = $finalPara.Range.End
$textRange = $d.Range($paraEnd - 1 - $oldText.Length, $paraEnd - 1)
$textRange.Text = $newText

Write-Output "Edit complete"
